# Data/Test Suite Prerequisites.xlsx
# - Add a new "patientNoAppointments" scenario row to the Patients sheet
# - Make the Patients sheet the active tab / selection
# (Shared-string table growth and the resulting index shift on the other
#  sheets happen automatically as a consequence of adding the two new
#  strings below - we don't touch those sheets directly.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Patients")

# Row 22 is the last existing data row; clone its formatting onto the new
# row 23 so the new row picks up the same borders/fills/alignment.
$ws.Range("B22:P22").Copy()
$ws.Range("B23:P23").PasteSpecial(-4122) # xlPasteFormats

# New scenario: a patient who should have no appointments at all.
$ws.Range("B23").Value = "patientNoAppointments"
$ws.Range("C23").Value = "?"
$ws.Range("D23").Value = "?"
$ws.Range("E23").Value = "?"
$ws.Range("F23").Value = "?"
$ws.Range("G23").Value = "?"
$ws.Range("H23").Value = "?"
$ws.Range("I23").Value = "?"
$ws.Range("J23").Value = "?"
$ws.Range("K23").Value = "?"
$ws.Range("L23").Value = "?"
$ws.Range("M23").Value = "?"
$ws.Range("N23").Value = "?"
$ws.Range("O23").Value = "?"
$ws.Range("P23").Value = "Patient should not have any appointments"

# Match the source row height (auto-grown to fit the wrapped label) and
# widen column B so "patientNoAppointments" is no longer clipped.
$ws.Rows.Item(23).RowHeight = 23.3
$ws.Columns.Item(2).ColumnWidth = 20.05

# The edited workbook re-opens on the Patients sheet with the new row
# selected (previously the Practitioner sheet was the active tab).
$ws.Activate()
$ws.Range("A23").Select()
